$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "FTR-left"
$ws.Range("D5").Value = "FTR-right"

$ws.Range("A4").Select()
